$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Row 2
$ws.Range("D2").Value = 0.32
$ws.Range("G2").Value = -0.006547251389746757
$ws.Range("H2").Value = -0.006547251389746757
$ws.Range("I2").Value = -0.0560840024706609
$ws.Range("J2").Value = -0.0560840024706609
$ws.Range("K2").Value = -7.56
$ws.Range("L2").Value = -0.04669549104385422
$ws.Range("U2").Value = 5.21
$ws.Range("V2").Value = 0.1069815195071868
$ws.Range("W2").Value = 2.719424460431655
$ws.Range("X2").Value = 0.1885631981122214
$ws.Range("Y2").Value = 2.530861262319433
$ws.Range("Z2").Value = 1.833729754219051
$ws.Range("AA2").Value = -0.1028429040661457
$ws.Range("AB2").Value = 0.1029841616526407
$ws.Range("AC2").Value = -0.2058270657187863
$ws.Range("AD2").Value = 70.40000000000001
$ws.Range("AF2").Value = 70.40000000000001
$ws.Range("AG2").Value = 65.19000000000001
$ws.Range("AH2").Value = 0.5910999160369438
$ws.Range("AI2").Value = 0.9777777777777779
$ws.Range("AJ2").Value = 0.5723944156642374
$ws.Range("AK2").Value = 0.9760443180116785
$ws.Range("AL2").Value = 8.32
$ws.Range("AM2").Value = 8.319000000000001
$ws.Range("AN2").Value = -32
$ws.Range("AO2").Value = -1.091346153846154
$ws.Range("AP2").Value = -29.63181818181819
$ws.Range("AQ2").Value = -1.091477341026566

# Row 3
$ws.Range("B3").Value = "AVE S.A. (ATSE:AVE)"
$ws.Range("D3").Value = 0.32
$ws.Range("G3").Value = -0.006547251389746757
$ws.Range("H3").Value = -0.006547251389746757
$ws.Range("I3").Value = -0.0560840024706609
$ws.Range("J3").Value = -0.0560840024706609
$ws.Range("K3").Value = -7.56
$ws.Range("L3").Value = -0.04669549104385422
$ws.Range("U3").Value = 5.21
$ws.Range("V3").Value = 0.1069815195071868
$ws.Range("W3").Value = 2.719424460431655
$ws.Range("X3").Value = 0.1885631981122214
$ws.Range("Y3").Value = 2.530861262319433
$ws.Range("Z3").Value = 1.833729754219051
$ws.Range("AA3").Value = -0.1028429040661457
$ws.Range("AB3").Value = 0.1029841616526407
$ws.Range("AC3").Value = -0.2058270657187863
$ws.Range("AD3").Value = 70.40000000000001
$ws.Range("AF3").Value = 70.40000000000001
$ws.Range("AG3").Value = 65.19000000000001
$ws.Range("AH3").Value = 0.5910999160369438
$ws.Range("AI3").Value = 0.9777777777777779
$ws.Range("AJ3").Value = 0.5723944156642374
$ws.Range("AK3").Value = 0.9760443180116785
$ws.Range("AL3").Value = 8.32
$ws.Range("AM3").Value = 8.319000000000001
$ws.Range("AN3").Value = -32
$ws.Range("AO3").Value = -1.091346153846154
$ws.Range("AP3").Value = -29.63181818181819
$ws.Range("AQ3").Value = -1.091477341026566
